$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.08900859163381439
$ws.Range("H2").Value = 34.21516522425406
$ws.Range("I2").Value = 37.57270272164795
$ws.Range("G3").Value = 0.135684081913601
$ws.Range("H3").Value = 14.73090971232015
$ws.Range("G4").Value = -0.3935677481740325
$ws.Range("H4").Value = -43.11796432013338
$ws.Range("G5").Value = -0.3338609746858758
$ws.Range("H5").Value = 16.32868724286925
$ws.Range("G6").Value = 0.2102466974330348
$ws.Range("H6").Value = 6.643454337288807
$ws.Range("G7").Value = 0.2741127412843961
$ws.Range("H7").Value = 32.17751379345965
$ws.Range("G8").Value = 0.1033787296731483
$ws.Range("H8").Value = 1.455122918893092
$ws.Range("G9").Value = 0.1325409579881081
$ws.Range("H9").Value = 4.791978032196555
$ws.Range("G10").Value = 0.05430473058539857
$ws.Range("H10").Value = -11.61023603449202
$ws.Range("G11").Value = 0.03393333233658978
$ws.Range("H11").Value = -32.03839445275428
$ws.Range("G12").Value = 0.1337172870606279
$ws.Range("H12").Value = 44.45426759433384
$ws.Range("G13").Value = 0.08785331497596693
$ws.Range("H13").Value = 15.28138687054351
$ws.Range("G14").Value = 0.2447859973650575
$ws.Range("H14").Value = 8.330774197711824
$ws.Range("G15").Value = 0.2538114800900541
$ws.Range("H15").Value = 3.023558034237123
$ws.Range("G16").Value = 0.1388227429153342
$ws.Range("H16").Value = 22.0474889182267
$ws.Range("G17").Value = 0.137494947518497
$ws.Range("H17").Value = -7.973740852697111
$ws.Range("G18").Value = 0.03369992356945477
$ws.Range("H18").Value = 476.453473568418
$ws.Range("G19").Value = 0.01561424089453524
$ws.Range("H19").Value = -35.53638521249197
$ws.Range("G20").Value = 0.112893995033526
$ws.Range("H20").Value = 32.72096896775692
$ws.Range("G21").Value = 0.0845498315713387
$ws.Range("H21").Value = 29.17449695003201
$ws.Range("G22").Value = 0.1712472278632153
$ws.Range("H22").Value = -10.60465277174532
$ws.Range("G23").Value = 0.2120633025769659
$ws.Range("H23").Value = -1.689335465230524
$ws.Range("G24").Value = 0.007709482702567048
$ws.Range("H24").Value = 302.5842507302712
$ws.Range("G25").Value = 0.01797141341366813
$ws.Range("H25").Value = 177.2720694616099
$ws.Range("G26").Value = 0.1913678783344877
$ws.Range("H26").Value = -6.58890291497327
$ws.Range("G27").Value = 0.2077986282000694
$ws.Range("H27").Value = 7.732132246610542
$ws.Range("G28").Value = 0.0261554003137069
$ws.Range("H28").Value = -60.91172536387048
$ws.Range("G29").Value = 0.1250790220639952
$ws.Range("H29").Value = 32.68931479349779
